$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2020-07-31 data: update nombre_aides (C) and montant_total (D) values.
# Values are stored as text in the sheet, so each new value is apostrophe-prefixed
# to force Excel to keep it as a text/string cell (matching the existing inlineStr cells)
# instead of auto-converting the numeric-looking text to a Number cell.
$ws.Range("C9").Value = "'48"
$ws.Range("D9").Value = "'117250.00"
$ws.Range("C10").Value = "'314"
$ws.Range("D10").Value = "'980466.74"
$ws.Range("C11").Value = "'130"
$ws.Range("D11").Value = "'522391.77"
$ws.Range("C12").Value = "'31"
$ws.Range("D12").Value = "'146500.00"
$ws.Range("C15").Value = "'89"
$ws.Range("D15").Value = "'212152.38"
$ws.Range("C16").Value = "'397"
$ws.Range("D16").Value = "'1173294.19"
$ws.Range("C17").Value = "'140"
$ws.Range("D17").Value = "'585458.00"
$ws.Range("C18").Value = "'41"
$ws.Range("D18").Value = "'186045.00"
$ws.Range("C19").Value = "'14"
$ws.Range("D19").Value = "'89716.00"
$ws.Range("C20").Value = "'14"
$ws.Range("D20").Value = "'30621.00"
$ws.Range("C33").Value = "'81"
$ws.Range("D33").Value = "'194640.00"
$ws.Range("C34").Value = "'465"
$ws.Range("D34").Value = "'1310793.53"
$ws.Range("C35").Value = "'185"
$ws.Range("D35").Value = "'816840.11"
$ws.Range("C36").Value = "'64"
$ws.Range("D36").Value = "'318974.00"
$ws.Range("C50").Value = "'85"
$ws.Range("D50").Value = "'235937.17"
$ws.Range("C51").Value = "'510"
$ws.Range("D51").Value = "'1598750.52"
$ws.Range("C52").Value = "'224"
$ws.Range("D52").Value = "'873868.65"
$ws.Range("C54").Value = "'23"
$ws.Range("D54").Value = "'125213.00"
$ws.Range("C56").Value = "'608"
$ws.Range("D56").Value = "'1461621.26"
$ws.Range("C57").Value = "'2992"
$ws.Range("D57").Value = "'8337567.73"
$ws.Range("C58").Value = "'1504"
$ws.Range("D58").Value = "'5299887.99"
$ws.Range("C59").Value = "'512"
$ws.Range("D59").Value = "'2207307.25"
$ws.Range("C60").Value = "'92"
$ws.Range("D60").Value = "'551041.00"
$ws.Range("C62").Value = "'236"
$ws.Range("D62").Value = "'540773.00"
$ws.Range("C68").Value = "'84"
$ws.Range("D68").Value = "'218542.41"
$ws.Range("C69").Value = "'362"
$ws.Range("D69").Value = "'1067946.54"
$ws.Range("C70").Value = "'141"
$ws.Range("D70").Value = "'544092.18"
$ws.Range("C71").Value = "'39"
$ws.Range("D71").Value = "'167497.67"
$ws.Range("C72").Value = "'9"
$ws.Range("D72").Value = "'51000.00"
$ws.Range("C74").Value = "'212"
$ws.Range("D74").Value = "'540326.09"
$ws.Range("C75").Value = "'820"
$ws.Range("D75").Value = "'2434044.56"
$ws.Range("C76").Value = "'302"
$ws.Range("D76").Value = "'1121266.79"
$ws.Range("C79").Value = "'28"
$ws.Range("D79").Value = "'57500.00"
$ws.Range("C86").Value = "'89"
$ws.Range("D86").Value = "'215878.00"
$ws.Range("C87").Value = "'382"
$ws.Range("D87").Value = "'1078760.67"
$ws.Range("C88").Value = "'163"
$ws.Range("D88").Value = "'620292.91"
$ws.Range("C90").Value = "'8"
$ws.Range("D90").Value = "'37500.00"
$ws.Range("C92").Value = "'248"
$ws.Range("D92").Value = "'618045.27"
